$d = $word.ActiveDocument

function Add-Paragraph($text) {
    $p = $d.Paragraphs.Last
    $r = $p.Range
    $r.Collapse(0)
    $r.InsertParagraphAfter()
    $newP = $d.Paragraphs.Last
    if ($text) {
        $newP.Range.InsertAfter($text)
    }
    return $newP
}


Add-Paragraph $null | Out-Null
Add-Paragraph "Flow of making ML Model" | Out-Null
Add-Paragraph "1. Import Libraries" | Out-Null
Add-Paragraph "2. Data Importing" | Out-Null
Add-Paragraph "3. Data Cleaning (Optional if Unitdy)" | Out-Null
Add-Paragraph "4. X & Y definition" | Out-Null
Add-Paragraph "5. Categorical _ cols" | Out-Null
Add-Paragraph "6. Data splitting into X_train, X_test, Y_train, Y_test" | Out-Null
Add-Paragraph "7. Standardization" | Out-Null
Add-Paragraph "8. Model selection" | Out-Null
Add-Paragraph "9. Model Training" | Out-Null
Add-Paragraph "10. then applying the model on the new test data" | Out-Null

Write-Output "Paragraphs: $($d.Paragraphs.Count)"
